# fee_waiver_template.docx edit:
#  1. Merge the split "{{ contact_details.suburb }}" / ".state" / ".postcode"
#     runs back into single runs (no visible text change).
#  2. Append a new "Comments to applicant:  {{ comments_to_applicant }}"
#     paragraph (plus two blank spacer paragraphs) right after the
#     "{% endfor %}" paragraph that closes the visits loop.

$d = $word.ActiveDocument

# --- 1. Run merges -------------------------------------------------------
# Re-"typing" the same text via Find/Replace collapses the two adjacent
# runs (identical formatting, no proofErr between them) into a single run,
# exactly like the diff shows.
$d.Content.Find.Execute("_details.suburb", $false, $false, $false, $false, $false, `
    $true, 1, $false, "_details.suburb", 2) | Out-Null
$d.Content.Find.Execute("_details.state", $false, $false, $false, $false, $false, `
    $true, 1, $false, "_details.state", 2) | Out-Null
$d.Content.Find.Execute("_details.postcode", $false, $false, $false, $false, $false, `
    $true, 1, $false, "_details.postcode", 2) | Out-Null

# --- 2. New "Comments to applicant" paragraphs ---------------------------
# Locate the paragraph that is exactly "{% endfor %}" (the one closing the
# visits loop, as opposed to the two inner "... {% endfor %}" occurrences).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Trim() -eq "{% endfor %}") {
        $target = $para
    }
}

# Create a placeholder paragraph right after it, then replace that
# placeholder's range with the real OOXML fragment (keeps the "{% endfor %}"
# paragraph itself untouched, unlike inserting straight into a collapsed
# range sitting on the paragraph boundary).
$target.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($target.Index + 1)
$insertionRange = $newPara.Range

$frag = '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Comments to applicant</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>:  {</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>comments_to_applicant</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> }}</w:t></w:r></w:p>'

$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    $frag + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionRange.InsertXML($xml)

Write-Host "Applied run merges + inserted 'Comments to applicant' paragraphs."
